$d = $word.ActiveDocument

# Replace a whole paragraph's text (excluding the trailing paragraph mark)
# by locating it via a unique substring marker. Using Range.Text directly
# (rather than Find.Execute) avoids Word's "smart quotes" AutoCorrect from
# mangling the straight ASCII quotes/backticks used in these lines.
function Set-ParagraphText([string]$marker, [string]$newText) {
    foreach ($p in $d.Paragraphs) {
        $r = $p.Range
        if ($r.Text.Contains($marker)) {
            $r.End = $r.End - 1
            $r.Text = $newText
            return
        }
    }
}

Set-ParagraphText 'New in this update (Railway Docker fix)' 'New in this update (Railway frontend Docker fix)'

Set-ParagraphText '- Fixed backend Dockerfile build context issue for Railway service root `backend-dotnet`.' '- Fixed frontend Dockerfile build context issue for Railway service root `frontend`.'

Set-ParagraphText '- Updated Dockerfile COPY commands to use local project context:' '- Updated Dockerfile COPY commands to use local context:'

Set-ParagraphText '  - `COPY *.csproj ./`' '  - `COPY package.json ./`'

Set-ParagraphText '- This resolves error: `"/backend-dotnet": not found`.' '- This resolves error: `"/frontend/package.json": not found`.'

Set-ParagraphText '- Last pushed commit: 8642aca' '- Last pushed commit: b90292a'

Set-ParagraphText '- Current Railway Dockerfile fix is local and not pushed yet.' '- Current frontend Dockerfile fix is local and not pushed yet.'
